# Apply updated crypto price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '26.835.87'
$ws.Range("E2").Value = '  +0.88%  '

# Row 3
$ws.Range("D3").Value = '1.646.88'
$ws.Range("E3").Value = '  +0.28%  '

# Row 4
$ws.Range("E4").Value = '  +0.26%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.28'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.85%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.501'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.37%  '

# Row 7
$ws.Range("E7").Value = '  +0.68%  '

# Row 8
$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.252'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.09%  '

# Row 9
$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0629'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.65%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.16'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.00%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0844'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.17%  '

# Row 12
$ws.Range("D12").Value = '1.870.72'
$ws.Range("E12").Value = '  -0.10%  '

# Row 13
$ws.Range("D13").Value = '1.658.75'
$ws.Range("E13").Value = '  +3.53%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.17'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.96%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.527'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.35%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.53'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.77%  '

# Row 17
$ws.Range("D17").Value = '26.833.86'
$ws.Range("E17").Value = '  +0.71%  '

# Row 18
$ws.Range("D18").Value = '0.0₃0737'
$ws.Range("E18").Value = '  -1.47%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '213.97'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.39%  '

# Row 20
$ws.Range("E20").Value = '  +0.80%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.38'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.09%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.47'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +16.08%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.26'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.63%  '

# Row 24
$ws.Range("E24").Value = '  -1.85%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.64'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.26%  '

# Row 26
$ws.Range("E26").Value = '  +0.59%  '

# Row 27
$ws.Range("E27").Value = '  -1.42%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.12'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.04%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.69'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.66%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0508'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.33%  '

# Row 31
$ws.Range("E31").Value = '  +0.68%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.34'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.07%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.00'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.14%  '

# Row 34
$ws.Range("D34").Value = '1.299.37'
$ws.Range("E34").Value = '  +2.31%  '

# Row 35
$ws.Range("E35").Value = '  -0.14%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.45'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.71%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0176'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.71%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.534'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.93%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.823'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.34%  '

# Row 40
$ws.Range("E40").Value = '  +0.85%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.811'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.09%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.23'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.33%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.34'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.20%  '

# Row 44
$ws.Range("D44").Value = '1.796.59'
$ws.Range("E44").Value = '  +0.73%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '61.89'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.30%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '91.60'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.34%  '

# Row 47
$ws.Range("E47").Value = '  +1.23%  '

# Row 48
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").Value = '0.0₆0102'
$ws.Range("E48").Value = '  -2.65%  '

# Row 49
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0524'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.41%  '

# Row 50
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.65'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.44%  '

# Row 51
$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0974'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.14%  '
